$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G2 and de-de!G2 share the same "Latest HO Xliff Generate Date" /
# "Correspond Handback DateTime" value for the first row - update both.
$wsOverview.Range("G2").Value = "2016-07-26 08:20:40"
$wsDeDe.Range("G2").Value = "2016-07-26 08:20:40"

# zh-cn row 2: Correspond Handoff Datetime (G2) and Correspond Handback DateTime (J2)
$wsZhCn.Range("G2").Value = "2016-07-26 08:20:30"
$wsZhCn.Range("J2").Value = "2016-07-26 08:21:15"

# de-de row 2: Correspond Handback DateTime (J2)
$wsDeDe.Range("J2").Value = "2016-07-26 08:21:29"
